# Daily attendance processing - 2026-01-30 18:18:05
#
# Normalises the order of names/emails listed in the "Recorded By" column
# (column G) of the Session Analysis Results sheet. A handful of rows had
# "System" recorded before the human/automation email that actually
# triggered the entry; this reorders those comma-separated lists into the
# canonical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact (case-sensitive) before -> after strings seen in the "Recorded By"
# column. Only an exact, full-value match is rewritten so already-correct
# rows (e.g. "System, backup@backdoor.com", "System" alone, or a lone
# email) are left untouched.
$map = @{}
$map["System, system, backup@backdoor.com"] = "System, backup@backdoor.com, system"
$map["System, admin@admin.com"]             = "admin@admin.com, System"
$map["System, dnasr281@gmail.com"]          = "dnasr281@gmail.com, System"
$map["dnasr281@gmail.com, admin@admin.com"] = "admin@admin.com, dnasr281@gmail.com"

$keys = @(
    "System, system, backup@backdoor.com",
    "System, admin@admin.com",
    "System, dnasr281@gmail.com",
    "dnasr281@gmail.com, admin@admin.com"
)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -ne $null) {
        foreach ($key in $keys) {
            if ($current.Equals($key)) {
                $cell.Value = $map[$key]
                break
            }
        }
    }
}
